# Auto-generated script to update FFXIV Leve profit calculation sheets
# with refreshed market-board price data (per scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 1727
$ws.Cells.Item(29, 9).Value = 734.3333
$ws.Cells.Item(29, 10).Value = 2152.4285
$ws.Cells.Item(29, 11).Value = 2202.9999
$ws.Cells.Item(29, 12).Value = 6457.2855
$ws.Cells.Item(29, 13).Value = -1921.9999
$ws.Cells.Item(29, 14).Value = -7019.2855
$ws.Cells.Item(43, 8).Value = 9260159
$ws.Cells.Item(43, 9).Value = 999.6667
$ws.Cells.Item(43, 11).Value = 999.6667
$ws.Cells.Item(43, 13).Value = -930.6667
$ws.Cells.Item(69, 8).Value = 3134.1667
$ws.Cells.Item(69, 10).Value = 3626.25
$ws.Cells.Item(69, 12).Value = 10878.75
$ws.Cells.Item(69, 14).Value = -12626.75
$ws.Cells.Item(70, 8).Value = 2855.2
$ws.Cells.Item(70, 9).Value = 2872.4443
$ws.Cells.Item(70, 11).Value = 8617.332900000001
$ws.Cells.Item(70, 13).Value = -8347.332900000001
$ws.Cells.Item(72, 8).Value = 3134.1667
$ws.Cells.Item(72, 10).Value = 3626.25
$ws.Cells.Item(72, 12).Value = 32636.25
$ws.Cells.Item(72, 14).Value = -41372.25
$ws.Cells.Item(73, 8).Value = 2855.2
$ws.Cells.Item(73, 9).Value = 2872.4443
$ws.Cells.Item(73, 11).Value = 8617.332900000001
$ws.Cells.Item(73, 13).Value = -7681.332900000001
$ws.Cells.Item(80, 8).Value = 1032.1666
$ws.Cells.Item(80, 10).Value = 773
$ws.Cells.Item(80, 12).Value = 2319
$ws.Cells.Item(80, 14).Value = -4315
$ws.Cells.Item(83, 8).Value = 1032.1666
$ws.Cells.Item(83, 10).Value = 773
$ws.Cells.Item(83, 12).Value = 6957
$ws.Cells.Item(83, 14).Value = -16941
$ws.Cells.Item(111, 8).Value = 2800.6365
$ws.Cells.Item(111, 9).Value = 1268.4546
$ws.Cells.Item(111, 10).Value = 4332.8184
$ws.Cells.Item(111, 11).Value = 3805.3638
$ws.Cells.Item(111, 12).Value = 12998.4552
$ws.Cells.Item(111, 13).Value = -738.3638000000001
$ws.Cells.Item(111, 14).Value = -19132.4552
$ws.Cells.Item(128, 8).Value = 37846
$ws.Cells.Item(128, 10).Value = 37846
$ws.Cells.Item(128, 12).Value = 37846
$ws.Cells.Item(128, 14).Value = -47806
$ws.Cells.Item(141, 8).Value = 6928.7
$ws.Cells.Item(141, 9).Value = 9547.083000000001
$ws.Cells.Item(141, 10).Value = 3001.125
$ws.Cells.Item(141, 11).Value = 28641.249
$ws.Cells.Item(141, 12).Value = 9003.375
$ws.Cells.Item(141, 13).Value = -23461.249
$ws.Cells.Item(141, 14).Value = -19363.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 1617.4
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 14).Value = -10900

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 2999.5
$ws.Cells.Item(11, 9).Value = 2999
$ws.Cells.Item(11, 10).Value = 3000
$ws.Cells.Item(11, 11).Value = 2999
$ws.Cells.Item(11, 12).Value = 3000
$ws.Cells.Item(11, 13).Value = -2859
$ws.Cells.Item(11, 14).Value = -3280
$ws.Cells.Item(18, 8).Value = 9000
$ws.Cells.Item(18, 10).Value = 9000
$ws.Cells.Item(18, 12).Value = 9000
$ws.Cells.Item(18, 14).Value = -10058
$ws.Cells.Item(20, 8).Value = 1216.2693
$ws.Cells.Item(20, 9).Value = 1100.7368
$ws.Cells.Item(20, 10).Value = 1529.8572
$ws.Cells.Item(20, 11).Value = 1100.7368
$ws.Cells.Item(20, 12).Value = 1529.8572
$ws.Cells.Item(20, 13).Value = -853.7367999999999
$ws.Cells.Item(20, 14).Value = -2023.8572
$ws.Cells.Item(86, 8).Value = 3359.6296
$ws.Cells.Item(86, 9).Value = 3584.8
$ws.Cells.Item(86, 10).Value = 2716.2856
$ws.Cells.Item(86, 11).Value = 3584.8
$ws.Cells.Item(86, 12).Value = 2716.2856
$ws.Cells.Item(86, 13).Value = -2461.8
$ws.Cells.Item(86, 14).Value = -4962.2856
$ws.Cells.Item(89, 8).Value = 3359.6296
$ws.Cells.Item(89, 9).Value = 3584.8
$ws.Cells.Item(89, 10).Value = 2716.2856
$ws.Cells.Item(89, 11).Value = 17924
$ws.Cells.Item(89, 12).Value = 13581.428
$ws.Cells.Item(89, 13).Value = -12308
$ws.Cells.Item(89, 14).Value = -24813.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1850.1364
$ws.Cells.Item(134, 9).Value = 1718.9375
$ws.Cells.Item(134, 10).Value = 2200
$ws.Cells.Item(134, 11).Value = 5156.8125
$ws.Cells.Item(134, 12).Value = 6600
$ws.Cells.Item(134, 13).Value = -2621.8125
$ws.Cells.Item(134, 14).Value = -11670

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1648.4517
$ws.Cells.Item(68, 10).Value = 2038.8572
$ws.Cells.Item(68, 12).Value = 6116.571599999999
$ws.Cells.Item(68, 14).Value = -7738.571599999999
$ws.Cells.Item(71, 8).Value = 1648.4517
$ws.Cells.Item(71, 10).Value = 2038.8572
$ws.Cells.Item(71, 12).Value = 18349.7148
$ws.Cells.Item(71, 14).Value = -26461.7148
$ws.Cells.Item(86, 8).Value = 500
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 500
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).ClearContents()
$ws.Cells.Item(86, 13).Value = 1500
$ws.Cells.Item(86, 14).Value = -3872
$ws.Cells.Item(89, 8).Value = 500
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 500
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).ClearContents()
$ws.Cells.Item(89, 13).Value = 4500
$ws.Cells.Item(89, 14).Value = -16356
$ws.Cells.Item(109, 8).Value = 74357.64
$ws.Cells.Item(109, 9).Value = 167867.83
$ws.Cells.Item(109, 10).Value = 4225
$ws.Cells.Item(109, 11).Value = 503603.49
$ws.Cells.Item(109, 12).Value = 12675
$ws.Cells.Item(109, 13).Value = -502563.49
$ws.Cells.Item(109, 14).Value = -14755

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 8).Value = 8000800
$ws.Cells.Item(24, 9).Value = 8000800
$ws.Cells.Item(24, 11).Value = 8000800
$ws.Cells.Item(24, 13).Value = -8000627
$ws.Cells.Item(80, 8).Value = 3329.9
$ws.Cells.Item(80, 9).Value = 1816.6666
$ws.Cells.Item(80, 11).Value = 1816.6666
$ws.Cells.Item(80, 13).Value = -818.6666
$ws.Cells.Item(83, 8).Value = 3329.9
$ws.Cells.Item(83, 9).Value = 1816.6666
$ws.Cells.Item(83, 11).Value = 9083.333000000001
$ws.Cells.Item(83, 13).Value = -4091.333000000001
$ws.Cells.Item(109, 8).Value = 30000
$ws.Cells.Item(109, 10).Value = 30000
$ws.Cells.Item(109, 12).Value = 30000
$ws.Cells.Item(109, 14).Value = -32080
$ws.Cells.Item(126, 8).Value = 2121.5
$ws.Cells.Item(126, 9).Value = 1737.375
$ws.Cells.Item(126, 10).Value = 2633.6667
$ws.Cells.Item(126, 11).Value = 5212.125
$ws.Cells.Item(126, 12).Value = 7901.000100000001
$ws.Cells.Item(126, 13).Value = -2742.125
$ws.Cells.Item(126, 14).Value = -12841.0001
$ws.Cells.Item(132, 8).Value = 2395.7646
$ws.Cells.Item(132, 9).Value = 2063.8076
$ws.Cells.Item(132, 11).Value = 6191.4228
$ws.Cells.Item(132, 13).Value = -3661.4228

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1144.3334
$ws.Cells.Item(22, 10).Value = 1285.5714
$ws.Cells.Item(22, 12).Value = 1285.5714
$ws.Cells.Item(22, 14).Value = -1875.5714
$ws.Cells.Item(27, 8).Value = 1144.3334
$ws.Cells.Item(27, 10).Value = 1285.5714
$ws.Cells.Item(27, 12).Value = 1285.5714
$ws.Cells.Item(27, 14).Value = -1499.5714
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 14).ClearContents()
$ws.Cells.Item(68, 8).Value = 1758
$ws.Cells.Item(68, 9).Value = 1302
$ws.Cells.Item(68, 11).Value = 1302
$ws.Cells.Item(68, 13).Value = -553
$ws.Cells.Item(71, 8).Value = 1758
$ws.Cells.Item(71, 9).Value = 1302
$ws.Cells.Item(71, 11).Value = 6510
$ws.Cells.Item(71, 13).Value = -2766
$ws.Cells.Item(82, 8).Value = 2815
$ws.Cells.Item(82, 9).Value = 2820
$ws.Cells.Item(82, 11).Value = 2820
$ws.Cells.Item(82, 13).Value = -2459
$ws.Cells.Item(85, 8).Value = 2815
$ws.Cells.Item(85, 9).Value = 2820
$ws.Cells.Item(85, 11).Value = 2820
$ws.Cells.Item(85, 13).Value = -1572

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(43, 8).Value = 10000
$ws.Cells.Item(43, 10).Value = 10000
$ws.Cells.Item(43, 12).Value = 10000
$ws.Cells.Item(43, 14).Value = -10298
$ws.Cells.Item(62, 8).Value = 83339030
$ws.Cells.Item(62, 9).Value = 100005840
$ws.Cells.Item(62, 10).Value = 5000
$ws.Cells.Item(62, 11).Value = 100005840
$ws.Cells.Item(62, 12).Value = 5000
$ws.Cells.Item(62, 13).Value = -100005216
$ws.Cells.Item(62, 14).Value = -6248
$ws.Cells.Item(65, 8).Value = 83339030
$ws.Cells.Item(65, 9).Value = 100005840
$ws.Cells.Item(65, 10).Value = 5000
$ws.Cells.Item(65, 11).Value = 500029200
$ws.Cells.Item(65, 12).Value = 25000
$ws.Cells.Item(65, 13).Value = -500026080
$ws.Cells.Item(65, 14).Value = -31240
$ws.Cells.Item(69, 8).Value = 26666.334
$ws.Cells.Item(69, 10).Value = 26666.334
$ws.Cells.Item(69, 12).Value = 26666.334
$ws.Cells.Item(69, 14).Value = -28164.334
$ws.Cells.Item(72, 8).Value = 26666.334
$ws.Cells.Item(72, 10).Value = 26666.334
$ws.Cells.Item(72, 12).Value = 79999.00199999999
$ws.Cells.Item(72, 14).Value = -87487.00199999999
$ws.Cells.Item(81, 8).Value = 1696.1875
$ws.Cells.Item(81, 9).Value = 1188.5
$ws.Cells.Item(81, 10).Value = 5250
$ws.Cells.Item(81, 11).Value = 2377
$ws.Cells.Item(81, 12).Value = 10500
$ws.Cells.Item(81, 13).Value = -1316
$ws.Cells.Item(81, 14).Value = -12622
$ws.Cells.Item(84, 8).Value = 1696.1875
$ws.Cells.Item(84, 9).Value = 1188.5
$ws.Cells.Item(84, 10).Value = 5250
$ws.Cells.Item(84, 11).Value = 11885
$ws.Cells.Item(84, 12).Value = 52500
$ws.Cells.Item(84, 13).Value = -6581
$ws.Cells.Item(84, 14).Value = -63108
$ws.Cells.Item(107, 8).Value = 591.1429000000001
$ws.Cells.Item(107, 9).Value = 472.25
$ws.Cells.Item(107, 11).Value = 1416.75
$ws.Cells.Item(107, 13).Value = 503.25
$ws.Cells.Item(126, 8).Value = 76924540
$ws.Cells.Item(126, 9).Value = 142857940
$ws.Cells.Item(126, 10).Value = 2249.1667
$ws.Cells.Item(126, 11).Value = 428573820
$ws.Cells.Item(126, 12).Value = 6747.500100000001
$ws.Cells.Item(126, 13).Value = -428571350
$ws.Cells.Item(126, 14).Value = -11687.5001
$ws.Cells.Item(136, 8).Value = 426.81482
$ws.Cells.Item(136, 9).Value = 359.33334
$ws.Cells.Item(136, 10).Value = 561.7778
$ws.Cells.Item(136, 11).Value = 1078.00002
$ws.Cells.Item(136, 12).Value = 1685.3334
$ws.Cells.Item(136, 13).Value = 1471.99998
$ws.Cells.Item(136, 14).Value = -6785.3334

